# Update generated stats on gh-pages output (南宁-漫展信息.xlsx)
$wb = $excel.ActiveWorkbook

# "展览" sheet: bump "想去人数" (F) counts for two events
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 7240
$ws1.Range("F5").Value = 178

# "全部类型" sheet mirrors the same records, update them too
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 7240
$ws4.Range("F6").Value = 178
